# Auto-sync update: insert three new user codes into the alphabetically
# sorted IDM_User_Codes list, shifting the existing rows below each
# insertion point down by one. The worksheet has a 1-row header
# ("username"/"code") starting at A1, with data starting at row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "KZU53477" / "1qwx-jHLR" before the old row 91 (LFA56301),
# i.e. right after KZL99357.
$ws.Rows("91:91").Insert()
$ws.Range("A91").Value = "KZU53477"
$ws.Range("B91").Value = "1qwx-jHLR"

# Insert "MXC93247" / "fcAI-rABo" before the old row 106/now-107
# (MXM72788), i.e. right after MUX13896.
$ws.Rows("107:107").Insert()
$ws.Range("A107").Value = "MXC93247"
$ws.Range("B107").Value = "fcAI-rABo"

# Insert "TKW11810" / "XbfA-yocA" before the old row 152/now-154
# (TMJ36970), i.e. right after SYH46751.
$ws.Rows("154:154").Insert()
$ws.Range("A154").Value = "TKW11810"
$ws.Range("B154").Value = "XbfA-yocA"
